# Refresh the cryptocurrency price/volume figures in columns D (Price) and
# E (Volume(1h)) on Sheet1, matching the GitHub Actions data-refresh commit.
#
# Some new Price values are digit-only strings (e.g. "23.20", "1.00") that
# Excel's automatic type detection would otherwise coerce into Numbers
# (dropping the significant trailing zero / reformatting the display).
# To keep them as literal text - exactly like the original inline-string
# cells - we briefly force NumberFormat "@" (Text) before assigning, then
# call ClearFormats() so the cell's style index reverts to the sheet's
# default (matching the unstyled source cells) without touching the text
# that was just written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.812.81'
$ws.Range("E2").Value = '  -0.47%  '

$ws.Range("D3").Value = '1.625.26'
$ws.Range("E3").Value = '  -0.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.996'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.87'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.41%  '

$ws.Range("E6").Value = '  -0.76%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.20'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.91%  '

$ws.Range("E9").Value = '  -1.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0607'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.00%  '

$ws.Range("E11").Value = '  -0.20%  '

$ws.Range("D12").Value = '1.856.24'

$ws.Range("D13").Value = '1.621.41'
$ws.Range("E13").Value = '  -0.42%  '

$ws.Range("E14").Value = '  -1.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.556'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.10%  '

$ws.Range("E16").Value = '  -0.86%  '

$ws.Range("D17").Value = '27.826.27'
$ws.Range("E17").Value = '  -0.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.97'
$ws.Range("D18").ClearFormats()

$ws.Range("D19").Value = '0.0₃0716'
$ws.Range("E19").Value = '  -1.16%  '

$ws.Range("E20").Value = '  +0.24%  '

$ws.Range("E21").Value = '  +0.01%  '

$ws.Range("E22").Value = '  -0.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.90'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.90%  '

$ws.Range("E24").Value = '  -0.71%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.93'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.36%  '

$ws.Range("E26").Value = '  -0.29%  '

$ws.Range("E27").Value = '  -0.45%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.41'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.995'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("E30").Value = '  -0.35%  '

$ws.Range("E31").Value = '  -0.39%  '

$ws.Range("E32").Value = '  -0.02%  '

$ws.Range("E33").Value = '  -0.15%  '

$ws.Range("D34").Value = '1.401.69'
$ws.Range("E34").Value = '  -0.21%  '

$ws.Range("E35").Value = '  +1.94%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.40%  '

$ws.Range("E37").Value = '  -1.16%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0169'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.93%  '

$ws.Range("E39").Value = '  -0.86%  '

$ws.Range("E40").Value = '  -3.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.996'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.08%  '

$ws.Range("E42").Value = '  -2.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.53'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.15%  '

$ws.Range("E44").Value = '  -1.69%  '

$ws.Range("E45").Value = '  -1.15%  '

$ws.Range("D46").Value = '1.765.59'
$ws.Range("E46").Value = '  -0.54%  '

$ws.Range("E47").Value = '  -3.85%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.07'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.14%  '

$ws.Range("E49").Value = '  +0.76%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0502'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.57'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.23%  '
